$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(5)

# --- Resize / reposition the shape ---
# Target EMU: off x=6744449 y=1536633 (unchanged); ext cx=1982861 cy=4799700 (unchanged)
# Shape.Left/.Top/.Width/.Height are points (Single); only Left & Width actually change.
$sh.Left = 531.0589913779527
$sh.Width = 156.13079140157478

# --- Update the notation text runs ---
$tr = $sh.TextFrame.TextRange

# Each paragraph holds "c.XXxx; e.XXxx" (or ";  " with two spaces). We replace the
# "c.XXxx" and "e.XXxx" tokens with underscore notation while leaving the
# "; "/";  " separator run untouched, splitting the single run into three.
$edits = @(
    @(2, "c_RPed", 10, "e_RPed"),
    @(18, "c_RPld", 27, "e_RPld"),
    @(39, "c_PCed", 47, "e_PCed"),
    @(55, "c_PCld", 64, "e_PCld"),
    @(76, "c_HCed", 84, "e_HCed"),
    @(92, "c_HCld", 101, "e_HCld")
)

foreach ($e in $edits) {
    $startC = $e[0]
    $newC = $e[1]
    $startE = $e[2]
    $newE = $e[3]

    $tr.Characters($startC, 6).Text = $newC
    $tr.Characters($startE, 6).Text = $newE
}
